# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.954.56'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '2.361.82'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'303.32"
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').Value = "'95.49"
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = "'0.483"
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('D10').Value = "'34.13"
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').Value = "'0.125"
$ws.Range('E11').Value = '  +3.33%  '
$ws.Range('D12').Value = "'0.0785"
$ws.Range('D13').Value = "'18.47"
$ws.Range('E13').Value = '  -3.32%  '
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '2.726.98'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').Value = '2.360.67'
$ws.Range('E16').Value = '  -0.81%  '
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').Value = '42.926.78'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = "'11.93"
$ws.Range('E19').Value = '  -2.45%  '
$ws.Range('D20').Value = "'6.26"
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').Value = '0.0₃0885'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('D22').Value = "'68.04"
$ws.Range('D23').Value = "'234.99"
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('E24').Value = '  -3.77%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = "'2.43"
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('E28').Value = '  +15.23%  '
$ws.Range('E29').Value = '  +1.80%  '
$ws.Range('D30').Value = "'32.08"
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').Value = "'5.00"
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').Value = "'17.53"
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('E35').Value = '  +2.94%  '
$ws.Range('E36').Value = '  +1.12%  '
$ws.Range('D37').Value = "'126.29"
$ws.Range('E37').Value = '  -11.24%  '
$ws.Range('E38').Value = '  -3.00%  '
$ws.Range('E39').Value = '  +2.62%  '
$ws.Range('D40').Value = "'2.25"
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('E42').Value = '  -3.73%  '
$ws.Range('D43').Value = '1.927.88'
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').Value = "'0.0278"
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('E45').Value = '  +3.94%  '
$ws.Range('D46').Value = "'9.21"
$ws.Range('E47').Value = '  -2.17%  '
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').Value = "'1.51"
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('D50').Value = "'71.50"
$ws.Range('E50').Value = '  -2.00%  '
$ws.Range('E51').Value = '  +0.97%  '
